$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data table: Dato / AntalTest / AntalOmikron / Ratio
# 21 rows (2021-11-27 .. 2021-12-17) replacing the previous 19-row table
# (which covered 2021-11-22 .. 2021-12-10), and filling in rows 21-22
# which previously only held the style with no data.

$dates = @(
    "2021-11-27","2021-11-28","2021-11-29","2021-11-30",
    "2021-12-01","2021-12-02","2021-12-03","2021-12-04","2021-12-05",
    "2021-12-06","2021-12-07","2021-12-08","2021-12-09","2021-12-10",
    "2021-12-11","2021-12-12","2021-12-13","2021-12-14","2021-12-15",
    "2021-12-16","2021-12-17"
)

$antalTest = @(
    4029,4000,5034,5374,4456,4563,5161,5301,5176,7115,7316,
    6575,6536,6823,6661,7048,9864,11220,10981,7605,3770
)

$antalOmikron = @(
    3,11,12,24,77,63,77,112,169,355,552,
    691,750,897,1113,1557,2863,4372,5095,3480,1632
)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $i + 2
    # Leading apostrophe keeps the date strings as literal text (matching
    # the existing quote-prefixed date-style cells in column A).
    $ws.Range("A$r").Value = "'" + $dates[$i]
    $ws.Range("B$r").Value = $antalTest[$i]
    $ws.Range("C$r").Value = $antalOmikron[$i]
}

# Ratio column: percentage of AntalOmikron over AntalTest
$ws.Range("D2:D22").Formula = "=100*C2/B2"

# Update the remembered selection shown in the workbook's last saved view
$ws.Range("B27").Select()
